$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 9).Value = 4273
$ws.Cells.Item(3, 9).Value = 4482
$ws.Cells.Item(4, 9).Value = 1032
$ws.Cells.Item(5, 9).Value = 408
$ws.Cells.Item(6, 9).Value = 4889
$ws.Cells.Item(7, 9).Value = 15084

$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(2, 9).Value = 43
$ws.Cells.Item(6, 9).Value = 63
$ws.Cells.Item(7, 9).Value = 173

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(2, 9).Value = 45
$ws.Cells.Item(4, 9).Value = 10
$ws.Cells.Item(7, 9).Value = 159

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 9).Value = 156
$ws.Cells.Item(3, 9).Value = 161
$ws.Cells.Item(6, 9).Value = 133
$ws.Cells.Item(7, 9).Value = 495

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(3, 9).Value = 99
$ws.Cells.Item(6, 9).Value = 103
$ws.Cells.Item(7, 9).Value = 339

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(4, 9).Value = 57
$ws.Cells.Item(7, 9).Value = 481
$ws.Cells.Item(8, 9).Value = 908
$ws.Cells.Item(9, 9).Value = 70
$ws.Cells.Item(10, 9).Value = 103
$ws.Cells.Item(11, 9).Value = 227
$ws.Cells.Item(16, 9).Value = 39
$ws.Cells.Item(17, 9).Value = 20
$ws.Cells.Item(19, 9).Value = 421
$ws.Cells.Item(20, 9).Value = 366
$ws.Cells.Item(23, 9).Value = 143
$ws.Cells.Item(24, 9).Value = 40
$ws.Cells.Item(25, 9).Value = 76
$ws.Cells.Item(27, 9).Value = 137
$ws.Cells.Item(28, 9).Value = 6
$ws.Cells.Item(29, 9).Value = 961
$ws.Cells.Item(32, 9).Value = 21
$ws.Cells.Item(33, 9).Value = 697
$ws.Cells.Item(34, 9).Value = 70
$ws.Cells.Item(37, 9).Value = 495
$ws.Cells.Item(42, 9).Value = 512
$ws.Cells.Item(43, 9).Value = 123
$ws.Cells.Item(44, 9).Value = 110
$ws.Cells.Item(47, 9).Value = 102
$ws.Cells.Item(49, 9).Value = 123
$ws.Cells.Item(52, 9).Value = 325
$ws.Cells.Item(54, 9).Value = 339
$ws.Cells.Item(55, 9).Value = 168
$ws.Cells.Item(57, 9).Value = 58
$ws.Cells.Item(60, 9).Value = 73
$ws.Cells.Item(63, 9).Value = 56
$ws.Cells.Item(65, 9).Value = 339
$ws.Cells.Item(68, 9).Value = 49
$ws.Cells.Item(72, 9).Value = 55
$ws.Cells.Item(73, 9).Value = 126
$ws.Cells.Item(76, 9).Value = 224
$ws.Cells.Item(78, 9).Value = 216
$ws.Cells.Item(79, 9).Value = 412
$ws.Cells.Item(80, 9).Value = 52
$ws.Cells.Item(82, 9).Value = 19
$ws.Cells.Item(83, 9).Value = 310
$ws.Cells.Item(85, 9).Value = 680
$ws.Cells.Item(86, 9).Value = 89
$ws.Cells.Item(87, 9).Value = 32
$ws.Cells.Item(88, 9).Value = 136
$ws.Cells.Item(89, 9).Value = 173
$ws.Cells.Item(95, 9).Value = 247
$ws.Cells.Item(96, 9).Value = 159
$ws.Cells.Item(101, 9).Value = 15084

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(6, 9).Value = 59
$ws.Cells.Item(7, 9).Value = 310

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(3, 9).Value = 96
$ws.Cells.Item(6, 9).Value = 43
$ws.Cells.Item(7, 9).Value = 247

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(2, 9).Value = 161
$ws.Cells.Item(6, 9).Value = 213
$ws.Cells.Item(7, 9).Value = 697

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Cells.Item(3, 9).Value = 14
$ws.Cells.Item(7, 9).Value = 123

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(2, 9).Value = 77
$ws.Cells.Item(3, 9).Value = 68
$ws.Cells.Item(7, 9).Value = 339

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 9).Value = 275
$ws.Cells.Item(6, 9).Value = 266
$ws.Cells.Item(7, 9).Value = 961

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(2, 9).Value = 153
$ws.Cells.Item(3, 9).Value = 123
$ws.Cells.Item(4, 9).Value = 18
$ws.Cells.Item(6, 9).Value = 118
$ws.Cells.Item(7, 9).Value = 421

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Cells.Item(2, 9).Value = 39
$ws.Cells.Item(7, 9).Value = 110

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(6, 9).Value = 96
$ws.Cells.Item(7, 9).Value = 224

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 9).Value = 175
$ws.Cells.Item(6, 9).Value = 170
$ws.Cells.Item(7, 9).Value = 680

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(2, 9).Value = 137
$ws.Cells.Item(4, 9).Value = 42
$ws.Cells.Item(7, 9).Value = 512

$ws = $wb.Worksheets.Item('Avondale')
$ws.Cells.Item(6, 9).Value = 45
$ws.Cells.Item(7, 9).Value = 103

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(5, 9).Value = 3
$ws.Cells.Item(7, 9).Value = 216

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(2, 9).Value = 52
$ws.Cells.Item(7, 9).Value = 168

$ws = $wb.Worksheets.Item('Dunning')
$ws.Cells.Item(2, 9).Value = 14
$ws.Cells.Item(7, 9).Value = 40

$ws = $wb.Worksheets.Item('Douglas')
$ws.Cells.Item(2, 9).Value = 40
$ws.Cells.Item(3, 9).Value = 49
$ws.Cells.Item(4, 9).Value = 9
$ws.Cells.Item(7, 9).Value = 143

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(3, 9).Value = 132
$ws.Cells.Item(7, 9).Value = 412

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(2, 9).Value = 100
$ws.Cells.Item(7, 9).Value = 366

$ws = $wb.Worksheets.Item('Burnside')
$ws.Cells.Item(3, 9).Value = 10
$ws.Cells.Item(7, 9).Value = 20

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(2, 9).Value = 92
$ws.Cells.Item(7, 9).Value = 325

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Cells.Item(3, 9).Value = 25
$ws.Cells.Item(7, 9).Value = 70

$ws = $wb.Worksheets.Item('East Side')
$ws.Cells.Item(3, 9).Value = 22
$ws.Cells.Item(7, 9).Value = 76

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Cells.Item(3, 9).Value = 33
$ws.Cells.Item(7, 9).Value = 102

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(2, 9).Value = 101
$ws.Cells.Item(7, 9).Value = 227

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Cells.Item(6, 9).Value = 16
$ws.Cells.Item(7, 9).Value = 70

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(2, 9).Value = 45
$ws.Cells.Item(6, 9).Value = 31
$ws.Cells.Item(7, 9).Value = 126

$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(6, 9).Value = 41
$ws.Cells.Item(7, 9).Value = 136

$ws = $wb.Worksheets.Item('Galewood')
$ws.Cells.Item(2, 9).Value = 7
$ws.Cells.Item(7, 9).Value = 21

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 9).Value = 286
$ws.Cells.Item(3, 9).Value = 260
$ws.Cells.Item(7, 9).Value = 908

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Cells.Item(3, 9).Value = 27
$ws.Cells.Item(6, 9).Value = 55
$ws.Cells.Item(7, 9).Value = 137

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Cells.Item(4, 9).Value = 44
$ws.Cells.Item(7, 9).Value = 89

$ws = $wb.Worksheets.Item('North Park')
$ws.Cells.Item(2, 9).Value = 19
$ws.Cells.Item(7, 9).Value = 49

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Cells.Item(6, 9).Value = 14
$ws.Cells.Item(7, 9).Value = 58

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Cells.Item(2, 9).Value = 20
$ws.Cells.Item(6, 9).Value = 22
$ws.Cells.Item(7, 9).Value = 73

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Cells.Item(6, 9).Value = 71
$ws.Cells.Item(7, 9).Value = 123

$ws = $wb.Worksheets.Item('Old Town')
$ws.Cells.Item(6, 9).Value = 31
$ws.Cells.Item(7, 9).Value = 55

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Cells.Item(3, 9).Value = 5
$ws.Cells.Item(6, 9).Value = 19

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Cells.Item(6, 9).Value = 28
$ws.Cells.Item(7, 9).Value = 52

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(3, 9).Value = 150
$ws.Cells.Item(6, 9).Value = 122
$ws.Cells.Item(7, 9).Value = 481

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Cells.Item(2, 9).Value = 22
$ws.Cells.Item(7, 9).Value = 57

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Cells.Item(6, 9).Value = 17
$ws.Cells.Item(7, 9).Value = 32

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Cells.Item(6, 9).Value = 24
$ws.Cells.Item(7, 9).Value = 39

$ws = $wb.Worksheets.Item('Edison Park')
$ws.Cells.Item(3, 9).Value = 3
$ws.Cells.Item(7, 9).Value = 6
